# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect the refreshed counts from the latest data pull.

$wb = $excel.ActiveWorkbook

# Row => new value for sheet "展览"
$exhibitionUpdates = @{
    5  = 153
    9  = 741
    11 = 64
    14 = 6403
    15 = 11
    20 = 15333
    21 = 1522
    22 = 281
    24 = 101
    25 = 11048
    26 = 750
    28 = 238
}

# Row => new value for sheet "全部类型"
$allTypesUpdates = @{
    5  = 153
    10 = 741
    13 = 64
    17 = 6403
    18 = 11
    23 = 15333
    24 = 1522
    25 = 281
    27 = 101
    28 = 11048
    29 = 750
    31 = 238
}

$wsExhibition = $wb.Worksheets.Item("展览")
foreach ($row in $exhibitionUpdates.Keys) {
    $wsExhibition.Range("F$row").Value = $exhibitionUpdates[$row]
}

$wsAllTypes = $wb.Worksheets.Item("全部类型")
foreach ($row in $allTypesUpdates.Keys) {
    $wsAllTypes.Range("F$row").Value = $allTypesUpdates[$row]
}
